$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 23:50"

# Country data rows (A:H), re-sorted by total cases descending, with updated stats
$data = @(
  @("Estados Unidos", 211408, 22878, 8805, 197885, 5005, 665, 4718),
  @("Italia", 110574, 4782, 16847, 80572, 4035, 727, 13155),
  @("España", 104118, 8195, 22647, 72084, 5872, 923, 9387),
  @("China", 81554, 36, 76238, 2004, 466, 7, 3312),
  @("Alemania", 77872, 6064, 18700, 58252, 3408, 145, 920),
  @("Francia", 56989, 4861, 10934, 42023, 6017, 509, 4032),
  @("Iran", 47593, 2988, 15473, 29084, 3871, 138, 3036),
  @("Reino Unido", 29474, 4324, 135, 26987, 163, 563, 2352),
  @("Suiza", 17768, 1163, 2967, 14313, 348, 55, 488),
  @("Turquia", 15679, 2148, 333, 15069, 847, 63, 277),
  @("Belgica", 13964, 1189, 2132, 11004, 1088, 123, 828),
  @("Paises Bajos", 13614, 1019, 250, 12191, 1053, 134, 1173),
  @("Austria", 10711, 531, 1436, 9129, 215, 18, 146),
  @("Corea del Sur", 9887, 101, 5567, 4155, 55, 3, 165),
  @("Canada", 9560, 948, 1615, 7833, 120, 11, 112),
  @("Portugal", 8251, 808, 43, 8021, 230, 27, 187),
  @("Brasil", 6836, 1119, 127, 6469, 296, 39, 240),
  @("Israel", 6092, 734, 241, 5825, 95, 6, 26),
  @("Suecia", 4947, 512, 103, 4605, 393, 59, 239),
  @("Australia", 4864, 101, 345, 4497, 50, 2, 22),
  @("Noruega", 4863, 222, 13, 4806, 105, 5, 44),
  @("Chequia", 3508, 200, 61, 3408, 70, 8, 39),
  @("Irlanda", 3447, 212, 5, 3357, 103, 14, 85),
  @("Dinamarca", 3107, 247, 894, 2109, 145, 14, 104),
  @("Chile", 3031, 293, 234, 2781, 31, 4, 16),
  @("Malasia", 2908, 142, 645, 2218, 102, 2, 45),
  @("Rusia", 2777, 440, 190, 2563, 8, 7, 24),
  @("Ecuador", 2748, 446, 58, 2597, 100, 14, 93),
  @("Polonia", 2554, 243, 47, 2464, 50, 10, 43),
  @("Rumania", 2460, 215, 252, 2116, 57, 10, 92),
  @("Luxemburgo", 2319, 141, 80, 2210, 31, 6, 29),
  @("Filipinas", 2311, 227, 50, 2165, 1, 8, 96),
  @("Japon", 2178, 0, 472, 1649, 69, 0, 57),
  @("Pakistan", 2118, 180, 94, 1997, 12, 1, 27),
  @("India", 1998, 601, 148, 1792, 0, 23, 58),
  @("Tailandia", 1771, 120, 505, 1254, 23, 2, 12),
  @("Arabia Saudita", 1720, 157, 264, 1440, 31, 6, 16),
  @("Indonesia", 1677, 149, 103, 1417, 0, 21, 157),
  @("Finlandia", 1446, 28, 10, 1419, 62, 0, 17),
  @("Grecia", 1415, 101, 52, 1313, 90, 1, 50),
  @("Sudafrica", 1380, 27, 50, 1325, 7, 0, 5),
  @("Peru", 1323, 258, 394, 891, 49, 8, 38),
  @("Republica Dominicana", 1284, 175, 9, 1218, 0, 6, 57),
  @("Islandia", 1220, 85, 225, 993, 12, 0, 2),
  @("Mexico", 1215, 121, 35, 1151, 1, 1, 29),
  @("Panama", 1181, 0, 9, 1142, 50, 0, 30),
  @("Colombia", 1065, 159, 39, 1009, 47, 1, 17),
  @("Serbia", 1060, 160, 42, 990, 62, 5, 28),
  @("Argentina", 1054, 0, 248, 778, 0, 1, 28),
  @("Singapur", 1000, 74, 245, 752, 24, 0, 3),
  @("Croacia", 963, 96, 73, 884, 34, 0, 6),
  @("Argelia", 847, 131, 61, 728, 0, 14, 58),
  @("Eslovenia", 841, 39, 10, 816, 31, 0, 15),
  @("Catar", 835, 54, 71, 762, 37, 0, 2),
  @("Emiratos Arabes Unidos", 814, 150, 61, 745, 2, 2, 8),
  @("Ucrania", 794, 149, 13, 761, 0, 3, 20),
  @("Estonia", 779, 34, 33, 741, 15, 1, 5),
  @("Egipto", 779, 69, 179, 548, 0, 6, 52),
  @("Hong Kong", 765, 50, 147, 614, 5, 0, 4),
  @("Irak", 728, 34, 182, 494, 0, 2, 52),
  @("Crucero", 712, 0, 603, 98, 15, 0, 11),
  @("Nueva Zelanda", 708, 61, 83, 624, 2, 0, 1),
  @("Marruecos", 654, 37, 29, 586, 1, 3, 39),
  @("Lituania", 581, 44, 7, 566, 11, 0, 8),
  @("Armenia", 571, 39, 31, 536, 30, 1, 4),
  @("Barein", 569, 2, 337, 228, 3, 0, 4),
  @("Hungria", 525, 33, 40, 465, 17, 4, 20),
  @("Libano", 479, 16, 43, 422, 5, 2, 14),
  @("Bosnia y Herzegovina", 459, 39, 19, 427, 1, 0, 13),
  @("Letonia", 446, 48, 1, 445, 3, 0, 0),
  @("Tunez", 423, 29, 5, 406, 10, 2, 12),
  @("Moldavia", 423, 70, 23, 395, 44, 1, 5),
  @("Bulgaria", 422, 23, 20, 392, 18, 2, 10),
  @("Eslovaquia", 400, 37, 3, 396, 1, 1, 1),
  @("Principado de Andorra", 390, 14, 10, 366, 12, 2, 14),
  @("Kazajistan", 380, 37, 26, 351, 6, 1, 3),
  @("Costa Rica", 375, 28, 4, 369, 9, 0, 2),
  @("Azerbaiyan", 359, 61, 26, 328, 7, 0, 5),
  @("Republica de Macedonia", 354, 25, 17, 326, 4, 2, 11),
  @("Uruguay", 338, 0, 41, 295, 12, 1, 2),
  @("Taiwan", 329, 7, 45, 279, 0, 0, 5),
  @("Republica de Chipre", 320, 58, 28, 283, 11, 1, 9),
  @("Kuwait", 317, 28, 80, 237, 14, 0, 0),
  @("Burkina Faso", 282, 21, 46, 220, 0, 2, 16),
  @("Reunion", 281, 34, 40, 241, 3, 0, 0),
  @("Jordania", 278, 4, 36, 237, 5, 0, 5),
  @("Albania", 259, 16, 67, 177, 7, 0, 15),
  @("Afganistan", 237, 63, 5, 228, 0, 0, 4),
  @("San Marino", 236, 0, 13, 197, 16, 0, 26),
  @("Camerun", 233, 40, 10, 217, 0, 0, 6),
  @("Vietnam", 218, 6, 63, 155, 3, 0, 0),
  @("Cuba", 212, 26, 12, 194, 3, 0, 6),
  @("Oman", 210, 18, 34, 175, 3, 0, 1),
  @("Ghana", 195, 34, 31, 159, 1, 0, 5),
  @("Costa de Marfil", 190, 11, 9, 180, 0, 0, 1),
  @("Senegal", 190, 15, 45, 144, 0, 1, 1),
  @("Malta", 188, 19, 2, 186, 2, 0, 0),
  @("Uzbekistan", 181, 9, 12, 167, 8, 0, 2),
  @("Nigeria", 174, 39, 9, 163, 0, 0, 2),
  @("Islas Feroe", 173, 4, 75, 98, 1, 0, 0),
  @("Honduras", 172, 31, 3, 159, 4, 3, 10),
  @("Bielorrusia", 163, 11, 53, 108, 2, 1, 2),
  @("Mauricio", 161, 18, 0, 155, 1, 1, 6),
  @("Sri Lanka", 146, 3, 21, 122, 5, 1, 3),
  @("Venezuela", 143, 8, 41, 99, 6, 0, 3),
  @("Martinica", 135, 7, 27, 105, 16, 0, 3),
  @("Estado de Palestina", 134, 15, 18, 115, 0, 0, 1),
  @("Brunei", 131, 2, 52, 78, 3, 0, 1),
  @("Guadalupe", 125, 11, 24, 95, 14, 2, 6),
  @("Montenegro", 123, 14, 0, 121, 4, 0, 2),
  @("Georgia", 117, 7, 23, 94, 6, 0, 0),
  @("Bolivia", 115, 8, 1, 107, 3, 1, 7),
  @("Kirguistan", 111, 4, 3, 108, 3, 0, 0),
  @("Consejo Danes para los Refugiados", 109, 11, 3, 97, 0, 1, 9),
  @("Camboya", 109, 0, 25, 84, 1, 0, 0),
  @("Mayotte", 94, 0, 10, 83, 3, 0, 1),
  @("Trinidad yTobago", 90, 3, 1, 84, 0, 2, 5),
  @("Ruanda", 82, 7, 0, 82, 0, 0, 0),
  @("Kenia", 81, 22, 3, 77, 2, 0, 1),
  @("Gibraltar", 81, 12, 34, 47, 0, 0, 0),
  @("Niger", 74, 40, 0, 69, 0, 2, 5),
  @("Paraguay", 69, 4, 1, 65, 4, 0, 3),
  @("Liechtenstein", 68, 0, 0, 68, 0, 0, 0),
  @("Isla de Man", 68, 8, 0, 67, 0, 1, 1),
  @("Madagascar", 57, 0, 0, 57, 6, 0, 0),
  @("Aruba", 55, 0, 1, 54, 0, 0, 0),
  @("Monaco", 55, 3, 2, 52, 2, 0, 1),
  @("Banglades", 54, 3, 25, 23, 1, 1, 6),
  @("Guayana Francesa", 51, 0, 15, 36, 0, 0, 0),
  @("Uganda", 44, 0, 0, 44, 0, 0, 0),
  @("Macao", 41, 0, 10, 31, 0, 0, 0),
  @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
  @("Guatemala", 39, 1, 12, 26, 1, 0, 1),
  @("Jamaica", 38, 0, 2, 33, 0, 1, 3),
  @("Polinesia Francesa", 37, 0, 0, 37, 1, 0, 0),
  @("Zambia", 36, 0, 0, 36, 0, 0, 0),
  @("Togo", 36, 2, 10, 24, 0, 1, 2),
  @("Barbados", 34, 0, 0, 34, 0, 0, 0),
  @("Republica de Yibuti", 33, 3, 0, 33, 0, 0, 0),
  @("El Salvador", 32, 0, 0, 31, 4, 0, 1),
  @("Guam", 32, 0, 0, 31, 0, 0, 1),
  @("Bermudas", 32, 0, 10, 22, 0, 0, 0),
  @("Mali", 31, 3, 0, 28, 0, 1, 3),
  @("Guinea", 30, 8, 0, 30, 0, 0, 0),
  @("Etiopia", 29, 3, 2, 27, 2, 0, 0),
  @("Islas Caimanes", 22, 8, 0, 21, 0, 0, 1),
  @("Bahamas", 21, 7, 1, 19, 0, 1, 1),
  @("Tanzania", 20, 1, 1, 18, 0, 0, 1),
  @("Congo", 19, 0, 0, 19, 0, 0, 0),
  @("Maldivas", 19, 1, 13, 6, 0, 0, 0),
  @("Gabon", 18, 2, 0, 17, 0, 0, 1),
  @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
  @("Nueva Caledonia", 16, 0, 0, 16, 0, 0, 0),
  @("Haiti", 16, 1, 1, 15, 0, 0, 0),
  @("San Martin (Parte Holandesa)", 16, 10, 6, 9, 0, 1, 1),
  @("Eritrea", 15, 0, 0, 15, 0, 0, 0),
  @("Guinea Ecuatorial", 15, 0, 1, 14, 0, 0, 0),
  @("Birmania", 15, 0, 0, 14, 0, 0, 1),
  @("San Martin (Parte Francesa)", 15, 0, 2, 12, 0, 0, 1),
  @("Namibia", 14, 3, 2, 12, 0, 0, 0),
  @("Mongolia", 14, 2, 2, 12, 0, 0, 0),
  @("Santa Lucia", 13, 0, 1, 12, 0, 0, 0),
  @("Benin", 13, 4, 1, 12, 0, 0, 0),
  @("Dominica", 12, 0, 0, 12, 0, 0, 0),
  @("Guyana", 12, 0, 0, 10, 0, 0, 2),
  @("Curazao", 11, 0, 3, 7, 0, 0, 1),
  @("Seychelles", 10, 0, 0, 10, 0, 0, 0),
  @("Libia", 10, 0, 0, 10, 0, 0, 0),
  @("Mozambique", 10, 2, 0, 10, 0, 0, 0),
  @("Laos", 10, 1, 0, 10, 0, 0, 0),
  @("Surinam", 10, 1, 0, 10, 0, 0, 0),
  @("Siria", 10, 0, 0, 8, 0, 0, 2),
  @("Groenlandia", 10, 0, 2, 8, 0, 0, 0),
  @("Suazilandia", 9, 0, 0, 9, 0, 0, 0),
  @("Guinea-Bisau", 9, 1, 0, 9, 0, 0, 0),
  @("Granada", 9, 0, 0, 9, 0, 0, 0),
  @("Montserrat", 9, 7, 0, 7, 0, 2, 2),
  @("San Cristobal y Nieves", 8, 0, 0, 8, 0, 0, 0),
  @("Zimbabue", 8, 0, 0, 7, 0, 0, 1),
  @("Angola", 8, 1, 1, 5, 0, 0, 2),
  @("Antigua y Barbuda", 7, 0, 0, 7, 0, 0, 0),
  @("Republica del Chad", 7, 0, 0, 7, 0, 0, 0),
  @("Sudan", 7, 0, 2, 3, 0, 0, 2),
  @("Santa Sede", 6, 0, 0, 6, 0, 0, 0),
  @("Islas Turcas y Caicos", 6, 1, 0, 6, 0, 0, 0),
  @("Liberia", 6, 3, 0, 6, 0, 0, 0),
  @("Cabo Verde", 6, 0, 0, 5, 0, 0, 1),
  @("San Bartolome", 6, 0, 1, 5, 0, 0, 0),
  @("Mauritania", 6, 0, 2, 3, 0, 0, 1),
  @("Fiyi", 5, 0, 0, 5, 0, 0, 0),
  @("Nepal", 5, 0, 1, 4, 0, 0, 0),
  @("Nicaragua", 5, 0, 0, 4, 0, 0, 1),
  @("Somalia", 5, 0, 1, 4, 0, 0, 0),
  @("Butan", 4, 0, 0, 4, 0, 0, 0),
  @("Botsuana", 4, 0, 0, 3, 0, 0, 1),
  @("Gambia", 4, 0, 2, 1, 0, 0, 1),
  @("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0),
  @("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
  @("Belice", 3, 0, 0, 3, 0, 0, 0),
  @("Burundi", 2, 0, 0, 2, 0, 0, 0),
  @("Anguila", 2, 0, 0, 2, 0, 0, 0),
  @("Bonaire, San Eustaquio y Saba", 2, 2, 0, 2, 0, 0, 0),
  @("Sierra Leona", 2, 1, 0, 2, 0, 0, 0),
  @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
  @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
  @("San Vicente y las Granadinas", 1, 0, 1, 0, 0, 0, 0)
)

$row = 4
foreach ($item in $data) {
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $item[$col - 1]
    }
    $row++
}